$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so values like "1.00" or
# "57.073.57" are stored as strings, matching the workbook's existing
# inline-string cell typing instead of being auto-coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "57.073.57"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "3.056.57"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "514.45"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("D6").Value = "139.90"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.431"
$ws.Range("E8").Value = "  -1.46%  "
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").Value = "0.107"
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("D12").Value = "3.588.94"
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("D14").Value = "25.35"
$ws.Range("E14").Value = "  -4.64%  "
$ws.Range("D15").Value = "0.0000161"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").Value = "57.177.04"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "3.061.02"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").Value = "12.93"
$ws.Range("E19").Value = "  -2.40%  "
$ws.Range("D20").Value = "8.02"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").Value = "331.99"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "0.497"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("E25").Value = "  +3.70%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").Value = "0.0₃0899"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("D28").Value = "6.30"
$ws.Range("E28").Value = "  -4.94%  "
$ws.Range("D29").Value = "7.09"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "1.81"
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("D31").Value = "20.73"
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("D32").Value = "1.15"
$ws.Range("E32").Value = "  -2.88%  "
$ws.Range("D33").Value = "154.74"
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("D34").Value = "27.07"
$ws.Range("E34").Value = "  +8.31%  "
$ws.Range("D35").Value = "4.43"
$ws.Range("E35").Value = "  -4.86%  "
$ws.Range("D36").Value = "5.81"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("D38").Value = "0.0668"
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("D39").Value = "3.097.83"
$ws.Range("E39").Value = "  +1.30%  "
$ws.Range("D40").Value = "36.84"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").Value = "3.86"
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "0.655"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").Value = "2.261.54"
$ws.Range("E44").Value = "  +3.16%  "
$ws.Range("D45").Value = "0.0257"
$ws.Range("E45").Value = "  +6.54%  "
$ws.Range("D46").Value = "1.36"
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "19.80"
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("D48").Value = "5.84"
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "0.919"
$ws.Range("E49").Value = "  -3.35%  "
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("D51").Value = "249.02"
$ws.Range("E51").Value = "  +6.81%  "

# Restore the default (General/Normal) style on column D so no stray
# cell-level style index is left behind on the price cells.
$ws.Range("D2:D51").Style = "Normal"
